$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.955.18"
$ws.Range("E2").Value = "  +0.31%  "
Set-TextValue "D3" "1.818.75"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "309.97"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.14%  "
Set-TextValue "D7" "0.4659"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.97%  "
Set-TextValue "D9" "0.07357"
$ws.Range("E9").Value = "  +0.03%  "
Set-TextValue "D10" "0.8723"
$ws.Range("E10").Value = "  -0.29%  "
Set-TextValue "D11" "20.28"
$ws.Range("E11").Value = "  -0.83%  "
Set-TextValue "D12" "1.804.11"
$ws.Range("E12").Value = "  +0.41%  "
Set-TextValue "D13" "5.419"
$ws.Range("E13").Value = "  +1.11%  "
Set-TextValue "D14" "0.07112"
$ws.Range("E14").Value = "  +0.94%  "
Set-TextValue "D15" "6.520"
$ws.Range("E15").Value = "  +0.18%  "
Set-TextValue "D16" "91.53"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +0.18%  "
Set-TextValue "D18" "0.000008709"
$ws.Range("E18").Value = "  +0.22%  "
Set-TextValue "D19" "1.003"
$ws.Range("E19").Value = "  +0.20%  "
Set-TextValue "D20" "14.66"
$ws.Range("E20").Value = "  -0.42%  "
Set-TextValue "D21" "26.980.31"
$ws.Range("E21").Value = "  +0.35%  "
Set-TextValue "D22" "5.295"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +0.44%  "
Set-TextValue "D24" "2.071.23"
$ws.Range("E24").Value = "  +2.89%  "
Set-TextValue "D25" "1.894"
$ws.Range("E25").Value = "  -0.42%  "
Set-TextValue "D26" "150.74"
$ws.Range("E26").Value = "  -0.58%  "
Set-TextValue "D27" "18.43"
$ws.Range("E27").Value = "  +0.33%  "
Set-TextValue "D28" "2.145"
$ws.Range("E28").Value = "  -0.11%  "
Set-TextValue "D29" "5.272"
$ws.Range("E29").Value = "  -0.70%  "
Set-TextValue "D30" "116.81"
$ws.Range("E30").Value = "  +0.77%  "
Set-TextValue "D31" "0.08897"
$ws.Range("E31").Value = "  -0.08%  "
Set-TextValue "D32" "0.7587"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +0.87%  "
Set-TextValue "D34" "4.507"
$ws.Range("E34").Value = "  +1.08%  "
Set-TextValue "D35" "2.910"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +0.19%  "
Set-TextValue "D37" "1.095"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue "D38" "0.05297"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  -0.81%  "
Set-TextValue "D40" "2.972"
$ws.Range("E40").Value = "  +1.82%  "
Set-TextValue "D41" "7.178"
$ws.Range("E41").Value = "  +0.36%  "
Set-TextValue "D42" "0.5294"
$ws.Range("E42").Value = "  -0.33%  "
Set-TextValue "D43" "2.335"
$ws.Range("E43").Value = "  -4.57%  "
Set-TextValue "D44" "0.1653"
$ws.Range("E44").Value = "  -0.42%  "
Set-TextValue "D45" "8.452"
$ws.Range("E45").Value = "  -0.10%  "
Set-TextValue "D46" "0.4873"
$ws.Range("E46").Value = "  -1.68%  "
Set-TextValue "D47" "10.42"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +0.18%  "
Set-TextValue "D49" "1.665"
$ws.Range("E49").Value = "  -0.41%  "
Set-TextValue "D50" "103.40"
$ws.Range("E50").Value = "  -0.07%  "
Set-TextValue "D51" "0.06295"
$ws.Range("E51").Value = "  +0.10%  "
